$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the Fitness column (C2:C252) with the corrected constant value
# produced by the fixed simulated annealing algorithm.
$ws.Range("C2:C252").Value = 7668
